$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("21 August 2019")
$ws.Range("C10").Value = "this topic is completed and updated on github"
